$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so that decimal-looking
# values (e.g. "302.20") are not silently coerced into numbers by Excel -
# the source data stores every Price/Volume cell as inline text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '43.814.88'
$ws.Range('E2').Value = '  -6.47%  '
$ws.Range('D3').Value = '2.589.91'
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '302.20'
$ws.Range('E5').Value = '  -2.16%  '
$ws.Range('D6').Value = '96.88'
$ws.Range('E6').Value = '  -4.50%  '
$ws.Range('E7').Value = '  -3.59%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.564'
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').Value = '36.89'
$ws.Range('E10').Value = '  -7.37%  '
$ws.Range('E11').Value = '  -3.90%  '
$ws.Range('E12').Value = '  -4.77%  '
$ws.Range('D13').Value = '2.988.40'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').Value = '2.591.00'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').Value = '0.899'
$ws.Range('E16').Value = '  -3.64%  '
$ws.Range('D17').Value = '14.43'
$ws.Range('E17').Value = '  -4.14%  '
$ws.Range('D18').Value = '43.774.76'
$ws.Range('E18').Value = '  -6.81%  '
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').Value = '0.0₃0986'
$ws.Range('E20').Value = '  -3.45%  '
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  -3.92%  '
$ws.Range('D22').Value = '73.72'
$ws.Range('E22').Value = '  +2.48%  '
$ws.Range('D23').Value = '266.76'
$ws.Range('E23').Value = '  -3.52%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '2.24'
$ws.Range('E24').Value = '  +2.73%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '2.95'
$ws.Range('E25').Value = '  -3.63%  '
$ws.Range('D26').Value = '29.37'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = '10.31'
$ws.Range('E28').Value = '  -3.61%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '38.29'
$ws.Range('E30').Value = '  -2.83%  '
$ws.Range('D31').Value = '6.20'
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('D32').Value = '3.60'
$ws.Range('E32').Value = '  -1.39%  '
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('D34').Value = '152.36'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  -3.66%  '
$ws.Range('D36').Value = '0.0821'
$ws.Range('E36').Value = '  -3.00%  '
$ws.Range('E37').Value = '  -5.16%  '
$ws.Range('E38').Value = '  -1.90%  '
$ws.Range('D39').Value = '24.32'
$ws.Range('E39').Value = '  +4.51%  '
$ws.Range('D40').Value = '17.13'
$ws.Range('E40').Value = '  +6.68%  '
$ws.Range('D41').Value = '3.62'
$ws.Range('E41').Value = '  -2.13%  '
$ws.Range('D42').Value = '0.0318'
$ws.Range('E42').Value = '  -5.16%  '
$ws.Range('E43').Value = '  -5.45%  '
$ws.Range('D44').Value = '2.041.98'
$ws.Range('E44').Value = '  -4.35%  '
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('D46').Value = '88.42'
$ws.Range('E46').Value = '  -6.27%  '
$ws.Range('D47').Value = '9.18'
$ws.Range('E47').Value = '  -4.25%  '
$ws.Range('D48').Value = '1.63'
$ws.Range('E48').Value = '  +5.23%  '
$ws.Range('D49').Value = '2.842.96'
$ws.Range('D50').Value = '106.16'
$ws.Range('E50').Value = '  -3.53%  '
$ws.Range('D51').Value = '0.193'
$ws.Range('E51').Value = '  -4.54%  '

# Restore the default style on the Price column so the workbook formatting
# matches the original (text is preserved via Excel general string storage).
$ws.Range("D2:D51").Style = "Normal"
